$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cell, $value)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") '27.074.80'
Set-TextValue $ws.Range("E2") '  -0.69%  '

Set-TextValue $ws.Range("D3") '1.825.59'
Set-TextValue $ws.Range("E3") '  +0.14%  '

Set-TextValue $ws.Range("D4") '1.002'
Set-TextValue $ws.Range("E4") '  -0.36%  '

Set-TextValue $ws.Range("D5") '312.85'
Set-TextValue $ws.Range("E5") '  -0.54%  '

Set-TextValue $ws.Range("D6") '1.002'
Set-TextValue $ws.Range("E6") '  -0.32%  '

Set-TextValue $ws.Range("D7") '0.4568'
Set-TextValue $ws.Range("E7") '  +6.74%  '

Set-TextValue $ws.Range("D8") '0.3723'
Set-TextValue $ws.Range("E8") '  +1.05%  '

Set-TextValue $ws.Range("D9") '0.07304'
Set-TextValue $ws.Range("E9") '  +0.87%  '

Set-TextValue $ws.Range("D10") '0.8600'
Set-TextValue $ws.Range("E10") '  -0.30%  '

Set-TextValue $ws.Range("D11") '20.88'
Set-TextValue $ws.Range("E11") '  -0.65%  '

Set-TextValue $ws.Range("D12") '1.819.00'
Set-TextValue $ws.Range("E12") '  -0.34%  '

Set-TextValue $ws.Range("D13") '6.695'
Set-TextValue $ws.Range("E13") '  +0.28%  '

Set-TextValue $ws.Range("D14") '93.02'
Set-TextValue $ws.Range("E14") '  +4.39%  '

Set-TextValue $ws.Range("D15") '5.354'
Set-TextValue $ws.Range("E15") '  +0.76%  '

Set-TextValue $ws.Range("D16") '0.07106'
Set-TextValue $ws.Range("E16") '  -0.05%  '

Set-TextValue $ws.Range("D17") '1.004'
Set-TextValue $ws.Range("E17") '  -0.36%  '

Set-TextValue $ws.Range("D18") '0.000008840'
Set-TextValue $ws.Range("E18") '  -0.44%  '

Set-TextValue $ws.Range("D19") '1.002'
Set-TextValue $ws.Range("E19") '  -0.28%  '

Set-TextValue $ws.Range("D20") '15.01'
Set-TextValue $ws.Range("E20") '  -0.38%  '

Set-TextValue $ws.Range("D21") '27.127.78'
Set-TextValue $ws.Range("E21") '  -0.56%  '

Set-TextValue $ws.Range("D22") '5.191'
Set-TextValue $ws.Range("E22") '  +0.81%  '

Set-TextValue $ws.Range("D23") '10.96'
Set-TextValue $ws.Range("E23") '  +0.71%  '

Set-TextValue $ws.Range("D24") '2.003'
Set-TextValue $ws.Range("E24") '  -0.19%  '

Set-TextValue $ws.Range("D25") '151.77'
Set-TextValue $ws.Range("E25") '  -1.12%  '

Set-TextValue $ws.Range("D26") '2.221'
Set-TextValue $ws.Range("E26") '  +4.44%  '

Set-TextValue $ws.Range("D27") '18.45'
Set-TextValue $ws.Range("E27") '  +0.27%  '

Set-TextValue $ws.Range("D28") '5.272'
Set-TextValue $ws.Range("E28") '  +0.61%  '

Set-TextValue $ws.Range("D29") '117.37'
Set-TextValue $ws.Range("E29") '  +0.86%  '

Set-TextValue $ws.Range("D30") '0.08881'
Set-TextValue $ws.Range("E30") '  -0.28%  '

Set-TextValue $ws.Range("D31") '1.192'
Set-TextValue $ws.Range("E31") '  -0.85%  '

Set-TextValue $ws.Range("D32") '0.7574'
Set-TextValue $ws.Range("E32") '  -0.40%  '

Set-TextValue $ws.Range("D33") '2.963'
Set-TextValue $ws.Range("E33") '  +5.35%  '

Set-TextValue $ws.Range("D34") '4.472'
Set-TextValue $ws.Range("E34") '  +0.18%  '

Set-TextValue $ws.Range("D35") '1.002'

Set-TextValue $ws.Range("D36") '1.102'
Set-TextValue $ws.Range("E36") '  -1.48%  '

Set-TextValue $ws.Range("D37") '0.01970'
Set-TextValue $ws.Range("E37") '  -0.12%  '

Set-TextValue $ws.Range("D38") '0.05283'
Set-TextValue $ws.Range("E38") '  -0.09%  '

Set-TextValue $ws.Range("D39") '0.5355'
Set-TextValue $ws.Range("E39") '  +6.13%  '

Set-TextValue $ws.Range("D40") '7.188'
Set-TextValue $ws.Range("E40") '  +0.69%  '

Set-TextValue $ws.Range("D41") '2.883'
Set-TextValue $ws.Range("E41") '  -0.97%  '

Set-TextValue $ws.Range("D42") '0.1717'
Set-TextValue $ws.Range("E42") '  +1.84%  '

Set-TextValue $ws.Range("D43") '0.5245'
Set-TextValue $ws.Range("E43") '  +10.20%  '

Set-TextValue $ws.Range("D44") '8.574'
Set-TextValue $ws.Range("E44") '  -0.78%  '

Set-TextValue $ws.Range("D45") '10.72'
Set-TextValue $ws.Range("E45") '  +1.13%  '

Set-TextValue $ws.Range("D46") '1.972'
Set-TextValue $ws.Range("E46") '  +8.71%  '

Set-TextValue $ws.Range("D47") '105.65'
Set-TextValue $ws.Range("E47") '  -1.02%  '

Set-TextValue $ws.Range("D48") '1.677'

Set-TextValue $ws.Range("E49") '  -0.37%  '

Set-TextValue $ws.Range("D50") '0.06416'
Set-TextValue $ws.Range("E50") '  +0.21%  '

Set-TextValue $ws.Range("B51") 'Aave'
Set-TextValue $ws.Range("C51") 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue $ws.Range("D51") '63.48'
Set-TextValue $ws.Range("E51") '  +0.56%  '
